$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42606.57230324074
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

$ws.Range("B4").Value = -30
$ws.Range("C4").Value = 66
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 85
$ws.Range("G4").Value = 11104
$ws.Range("H4").Value = 9033
$ws.Range("I4").Value = 1515
$ws.Range("J4").Value = 184
$ws.Range("K4").Value = 90
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 53
$ws.Range("N4").Value = "Noun"
